$wb = $excel.ActiveWorkbook

# ---- ALC (39 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3522.5
$ws.Range("I64").Value = 3040
$ws.Range("J64").Value = 4326.6665
$ws.Range("K64").Value = 3040
$ws.Range("L64").Value = 4326.6665
$ws.Range("M64").Value = -2792
$ws.Range("N64").Value = -4822.6665
$ws.Range("H67").Value = 3522.5
$ws.Range("I67").Value = 3040
$ws.Range("J67").Value = 4326.6665
$ws.Range("K67").Value = 3040
$ws.Range("L67").Value = 4326.6665
$ws.Range("M67").Value = -2182
$ws.Range("N67").Value = -6042.6665
$ws.Range("H101").Value = 1278.8462
$ws.Range("I101").Value = 637.6667
$ws.Range("J101").Value = 1828.4286
$ws.Range("K101").Value = 1913.0001
$ws.Range("L101").Value = 5485.2858
$ws.Range("M101").Value = -291.0001
$ws.Range("N101").Value = -8729.2858
$ws.Range("H132").Value = 30914312
$ws.Range("I132").Value = 34832930
$ws.Range("K132").Value = 104498790
$ws.Range("M132").Value = -104496260
$ws.Range("H137").Value = 759065.9
$ws.Range("I137").Value = 2168387.8
$ws.Range("J137").Value = 2844.4146
$ws.Range("K137").Value = 6505163.399999999
$ws.Range("L137").Value = 8533.2438
$ws.Range("M137").Value = -6502613.399999999
$ws.Range("N137").Value = -13633.2438
$ws.Range("H138").Value = 2728.652
$ws.Range("I138").Value = 1911.3572
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 5734.071599999999
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -594.0715999999993
$ws.Range("N138").Value = -22280

# ---- ARM (44 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3851.4824
$ws.Range("I32").Value = 3715.757
$ws.Range("J32").Value = 4484.8667
$ws.Range("K32").Value = 3715.757
$ws.Range("L32").Value = 4484.8667
$ws.Range("M32").Value = -3428.757
$ws.Range("N32").Value = -5058.8667
$ws.Range("H61").Value = 4003
$ws.Range("I61").Value = 4670.6665
$ws.Range("K61").Value = 4670.6665
$ws.Range("M61").Value = -4458.6665
$ws.Range("H74").Value = 411648.12
$ws.Range("I74").Value = 969833.4399999999
$ws.Range("J74").Value = 2312.2
$ws.Range("K74").Value = 969833.4399999999
$ws.Range("L74").Value = 2312.2
$ws.Range("M74").Value = -968959.4399999999
$ws.Range("N74").Value = -4060.2
$ws.Range("H77").Value = 411648.12
$ws.Range("I77").Value = 969833.4399999999
$ws.Range("J77").Value = 2312.2
$ws.Range("K77").Value = 4849167.199999999
$ws.Range("L77").Value = 11561
$ws.Range("M77").Value = -4844799.199999999
$ws.Range("N77").Value = -20297
$ws.Range("H109").Value = 30122.125
$ws.Range("J109").Value = 30122.125
$ws.Range("L109").Value = 30122.125
$ws.Range("N109").Value = -32896.125
$ws.Range("H119").Value = 32676.309
$ws.Range("J119").Value = 32676.309
$ws.Range("L119").Value = 32676.309
$ws.Range("N119").Value = -42352.309
$ws.Range("H132").Value = 2542.2964
$ws.Range("I132").Value = 2079.9756
$ws.Range("J132").Value = 4000.3845
$ws.Range("K132").Value = 6239.926800000001
$ws.Range("L132").Value = 12001.1535
$ws.Range("M132").Value = -3709.926800000001
$ws.Range("N132").Value = -17061.1535
$ws.Range("H136").Value = 4003
$ws.Range("I136").Value = 4670.6665
$ws.Range("K136").Value = 14011.9995
$ws.Range("M136").Value = -11461.9995

# ---- BSM (27 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 44500
$ws.Range("J74").Value = 44500
$ws.Range("L74").Value = 44500
$ws.Range("N74").Value = -46372
$ws.Range("H77").Value = 44500
$ws.Range("J77").Value = 44500
$ws.Range("L77").Value = 133500
$ws.Range("N77").Value = -142860
$ws.Range("H103").Value = 34750
$ws.Range("J103").Value = 34750
$ws.Range("L103").Value = 34750
$ws.Range("N103").Value = -37094
$ws.Range("H125").Value = 41798.89
$ws.Range("J125").Value = 41798.89
$ws.Range("L125").Value = 41798.89
$ws.Range("N125").Value = -51638.89
$ws.Range("H134").Value = 3849.7632
$ws.Range("I134").Value = 1380.3334
$ws.Range("J134").Value = 5460.2607
$ws.Range("K134").Value = 4141.0002
$ws.Range("L134").Value = 16380.7821
$ws.Range("M134").Value = -1606.0002
$ws.Range("N134").Value = -21450.7821
$ws.Range("H139").Value = 64603.5
$ws.Range("I139").Value = 30707
$ws.Range("K139").Value = 30707
$ws.Range("M139").Value = -25567

# ---- CRP (58 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 378676.97
$ws.Range("I31").Value = 3380703.2
$ws.Range("J31").Value = 3423.6875
$ws.Range("K31").Value = 3380703.2
$ws.Range("L31").Value = 3423.6875
$ws.Range("M31").Value = -3380408.2
$ws.Range("N31").Value = -4013.6875
$ws.Range("H34").Value = 378676.97
$ws.Range("I34").Value = 3380703.2
$ws.Range("J34").Value = 3423.6875
$ws.Range("K34").Value = 3380703.2
$ws.Range("L34").Value = 3423.6875
$ws.Range("M34").Value = -3380501.2
$ws.Range("N34").Value = -3827.6875
$ws.Range("H62").Value = 2990.2727
$ws.Range("I62").Value = 2789.3
$ws.Range("K62").Value = 2789.3
$ws.Range("M62").Value = -2165.3
$ws.Range("H65").Value = 2990.2727
$ws.Range("I65").Value = 2789.3
$ws.Range("K65").Value = 13946.5
$ws.Range("M65").Value = -10826.5
$ws.Range("H87").Value = 25431.666
$ws.Range("J87").Value = 25431.666
$ws.Range("L87").Value = 25431.666
$ws.Range("N87").Value = -27803.666
$ws.Range("H90").Value = 25431.666
$ws.Range("J90").Value = 25431.666
$ws.Range("L90").Value = 76294.99800000001
$ws.Range("N90").Value = -88150.99800000001
$ws.Range("H99").Value = 3477.7273
$ws.Range("I99").Value = 1600.5385
$ws.Range("J99").Value = 6189.222
$ws.Range("K99").Value = 1600.5385
$ws.Range("L99").Value = 6189.222
$ws.Range("M99").Value = -102.5385000000001
$ws.Range("N99").Value = -9185.222
$ws.Range("H103").Value = 31887.334
$ws.Range("I103").Value = 19108
$ws.Range("J103").Value = 44666.668
$ws.Range("K103").Value = 19108
$ws.Range("L103").Value = 44666.668
$ws.Range("M103").Value = -17936
$ws.Range("N103").Value = -47010.668
$ws.Range("H126").Value = 3477.7273
$ws.Range("I126").Value = 1600.5385
$ws.Range("J126").Value = 6189.222
$ws.Range("K126").Value = 4801.6155
$ws.Range("L126").Value = 18567.666
$ws.Range("M126").Value = -2331.6155
$ws.Range("N126").Value = -23507.666
$ws.Range("H134").Value = 7606.0586
$ws.Range("I134").Value = 7548.933
$ws.Range("J134").Value = 8034.5
$ws.Range("K134").Value = 22646.799
$ws.Range("L134").Value = 24103.5
$ws.Range("M134").Value = -20111.799
$ws.Range("N134").Value = -29173.5

# ---- CUL (36 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1723073.1
$ws.Range("I4").Value = 12054762
$ws.Range("K4").Value = 36164286
$ws.Range("M4").Value = -36164174
$ws.Range("H68").Value = 3642.9111
$ws.Range("I68").Value = 1084.08
$ws.Range("J68").Value = 6841.45
$ws.Range("K68").Value = 3252.24
$ws.Range("L68").Value = 20524.35
$ws.Range("M68").Value = -2441.24
$ws.Range("N68").Value = -22146.35
$ws.Range("H71").Value = 3642.9111
$ws.Range("I71").Value = 1084.08
$ws.Range("J71").Value = 6841.45
$ws.Range("K71").Value = 9756.719999999999
$ws.Range("L71").Value = 61573.05
$ws.Range("M71").Value = -5700.719999999999
$ws.Range("N71").Value = -69685.04999999999
$ws.Range("H107").Value = 14134.368
$ws.Range("I107").Value = 430.45456
$ws.Range("J107").Value = 24651.326
$ws.Range("K107").Value = 1291.36368
$ws.Range("L107").Value = 73953.978
$ws.Range("M107").Value = 628.6363200000001
$ws.Range("N107").Value = -77793.978
$ws.Range("H131").Value = 798.79
$ws.Range("I131").Value = 454.4
$ws.Range("J131").Value = 837.05554
$ws.Range("K131").Value = 1363.2
$ws.Range("L131").Value = 2511.16662
$ws.Range("M131").Value = 3676.8
$ws.Range("N131").Value = -12591.16662
$ws.Range("H132").Value = 2572.5
$ws.Range("I132").Value = 1039.375
$ws.Range("K132").Value = 9354.375
$ws.Range("M132").Value = -6824.375

# ---- GSM (22 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 34000
$ws.Range("J4").Value = 34000
$ws.Range("L4").Value = 34000
$ws.Range("N4").Value = -34224
$ws.Range("H97").Value = 805
$ws.Range("I97").Value = 805
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 805
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -309
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 10101756
$ws.Range("J107").Value = 27779230
$ws.Range("L107").Value = 27779230
$ws.Range("N107").Value = -27783070
$ws.Range("H132").Value = 2692.0408
$ws.Range("I132").Value = 2048.3901
$ws.Range("J132").Value = 5990.75
$ws.Range("K132").Value = 6145.1703
$ws.Range("L132").Value = 17972.25
$ws.Range("M132").Value = -3615.1703
$ws.Range("N132").Value = -23032.25

# ---- LTW (11 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6958.8
$ws.Range("I122").Value = 3931.3333
$ws.Range("K122").Value = 11793.9999
$ws.Range("M122").Value = -9343.999899999999
$ws.Range("H136").Value = 5010.5835
$ws.Range("I136").Value = 2840.818
$ws.Range("J136").Value = 6846.5386
$ws.Range("K136").Value = 8522.454000000002
$ws.Range("L136").Value = 20539.6158
$ws.Range("M136").Value = -5972.454000000002
$ws.Range("N136").Value = -25639.6158

# ---- WVR (14 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 696.5833
$ws.Range("I107").Value = 520.5625
$ws.Range("J107").Value = 1048.625
$ws.Range("K107").Value = 1561.6875
$ws.Range("L107").Value = 3145.875
$ws.Range("M107").Value = 358.3125
$ws.Range("N107").Value = -6985.875
$ws.Range("H132").Value = 2073.6511
$ws.Range("I132").Value = 1241.7142
$ws.Range("J132").Value = 3626.6
$ws.Range("K132").Value = 3725.1426
$ws.Range("L132").Value = 10879.8
$ws.Range("M132").Value = -1195.1426
$ws.Range("N132").Value = -15939.8

Write-Host "Applied 251 cell updates across 8 sheets."